$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap Sender Email / Sender Name values in row 2 (A2 <-> B2)
$ws.Range("A2").Value = "mongodbteam@mongodb.com"
$ws.Range("B2").Value = "MongoDB"

# Remove the now-unused Relevant/Comment values on row 2
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

# Delete the duplicate data row entirely (row 3)
$ws.Rows.Item(3).Delete()
